$wb = $excel.ActiveWorkbook

# "Repayment Schedule" sheet: insert a new blank column at N, shifting the
# old N/O/P (Late / blank / Outstanding) columns right to O/P/Q.
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Columns("N").Insert() | Out-Null

# Give the freshly inserted column a fixed width (matches its neighbour M,
# minus the "best fit" auto-size flag that the other data columns carry).
$ws.Columns("N").ColumnWidth = 10.3

# Make "Repayment Schedule" the active sheet/tab (was "Transactions"),
# and move its selection to S8.
$ws.Activate() | Out-Null
$ws.Range("S8").Select() | Out-Null
